$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 18

$ws.Cells.Item($row, 1).Value = "MF464269"
$ws.Cells.Item($row, 2).Value = "Marmot-BV"
$ws.Cells.Item($row, 3).Value = "Himalayan marmot bocaparvovirus 2"
$ws.Cells.Item($row, 4).Value = "Marmota himalayana"
$ws.Cells.Item($row, 5).Value = "NK"
$ws.Cells.Item($row, 6).Value = "Parvovirinae"
$ws.Cells.Item($row, 7).Value = "Bocaparvovirus"
$ws.Cells.Item($row, 8).Value = "NK"
$ws.Cells.Item($row, 9).Value = "NK"
$ws.Cells.Item($row, 10).Value = "NK"
$ws.Cells.Item($row, 11).Value = "NK"
$ws.Cells.Item($row, 12).Value = "NK"

$ws.Range("A18").Select()
